# bill sheet: reconcile bill_status / bill_type values, drop two stale
# voucher-id (vou_id) entries, and refresh the sheet view/column widths
# to match the user's last interaction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: bill_status "Da hoan" -> "Xac nhan"; bill_type "Dong coc" -> "Hoan coc"
$ws.Range("D3").Value = "Xac nhan"
$ws.Range("E3").Value = "Hoan coc"

# Row 4: bill_status "Da hoan" -> "Xac nhan"; bill_type "Dong coc" -> "Hoan coc"; vou_id cleared
$ws.Range("D4").Value = "Xac nhan"
$ws.Range("E4").Value = "Hoan coc"
$ws.Range("G4").ClearContents()

# Row 8: bill_status "Da hoan" -> "Xac nhan"; vou_id cleared
$ws.Range("D8").Value = "Xac nhan"
$ws.Range("E8").Value = "Dong coc"
$ws.Range("G8").ClearContents()

# Row 9: bill_status unchanged text ("Chua xac nhan")
$ws.Range("D9").Value = "Chua xac nhan"

# Row 10: bill_status "Chua hoan" -> "Chua xac nhan"; bill_type "Dong coc" -> "Hoan coc"
$ws.Range("D10").Value = "Chua xac nhan"
$ws.Range("E10").Value = "Hoan coc"

# Column widths nudged slightly wider after the edits (bill_status / bill_type)
$ws.Columns("D").ColumnWidth = 13.5
$ws.Columns("E").ColumnWidth = 11.833333333333334

# Final view state: zoomed to 102%, last cell selected is G8
$ws.Range("G8").Select()
$excel.ActiveWindow.Zoom = 102
